$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style source cells for id-column (A) and date-column (E), copied from row 198
$idStyleSrc = $ws.Range("A198")
$dateStyleSrc = $ws.Range("E198")

# --- Update existing rows 199-202: add new cells / change values ---
# Row 199
$ws.Cells.Item(199, 8).Value = 4
$ws.Cells.Item(199, 9).Value = 0
$ws.Cells.Item(199, 10).Value = "H"
$ws.Cells.Item(199, 23).Value = 0.6499999999999999
$ws.Cells.Item(199, 24).Value = -1
$ws.Cells.Item(199, 25).Value = -1
$ws.Cells.Item(199, 26).Value = 0.8500000000000001
$ws.Cells.Item(199, 27).Value = -1
$ws.Cells.Item(199, 28).Value = 1.025
$ws.Cells.Item(199, 29).Value = -1

# Row 200
$ws.Cells.Item(200, 8).Value = 1
$ws.Cells.Item(200, 9).Value = 1
$ws.Cells.Item(200, 10).Value = "D"
$ws.Cells.Item(200, 14).Value = 2.875
$ws.Cells.Item(200, 15).Value = 3.3
$ws.Cells.Item(200, 16).Value = 2.45
$ws.Cells.Item(200, 17).Value = 0
$ws.Cells.Item(200, 18).Value = 2.1
$ws.Cells.Item(200, 19).Value = 1.775
$ws.Cells.Item(200, 21).Value = 1.95
$ws.Cells.Item(200, 22).Value = 1.9
$ws.Cells.Item(200, 23).Value = -1
$ws.Cells.Item(200, 24).Value = 2.3
$ws.Cells.Item(200, 25).Value = -1
$ws.Cells.Item(200, 28).Value = -0.5
$ws.Cells.Item(200, 29).Value = 0.45

# Row 201
$ws.Cells.Item(201, 8).Value = 3
$ws.Cells.Item(201, 9).Value = 0
$ws.Cells.Item(201, 10).Value = "H"
$ws.Cells.Item(201, 14).Value = 1.444
$ws.Cells.Item(201, 15).Value = 4.333
$ws.Cells.Item(201, 17).Value = -1.25
$ws.Cells.Item(201, 18).Value = 2
$ws.Cells.Item(201, 19).Value = 1.85
$ws.Cells.Item(201, 20).Value = 2.75
$ws.Cells.Item(201, 21).Value = 2.025
$ws.Cells.Item(201, 22).Value = 1.825
$ws.Cells.Item(201, 23).Value = 0.444
$ws.Cells.Item(201, 24).Value = -1
$ws.Cells.Item(201, 25).Value = -1
$ws.Cells.Item(201, 26).Value = 1
$ws.Cells.Item(201, 27).Value = -1
$ws.Cells.Item(201, 28).Value = 0.5125
$ws.Cells.Item(201, 29).Value = -0.5

# Row 202
$ws.Cells.Item(202, 8).Value = 2
$ws.Cells.Item(202, 9).Value = 0
$ws.Cells.Item(202, 10).Value = "H"
$ws.Cells.Item(202, 14).Value = 1.222
$ws.Cells.Item(202, 15).Value = 7
$ws.Cells.Item(202, 18).Value = 1.9
$ws.Cells.Item(202, 19).Value = 1.95
$ws.Cells.Item(202, 21).Value = 1.95
$ws.Cells.Item(202, 22).Value = 1.9
$ws.Cells.Item(202, 23).Value = 0.222
$ws.Cells.Item(202, 24).Value = -1
$ws.Cells.Item(202, 25).Value = -1
$ws.Cells.Item(202, 26).Value = 0.45
$ws.Cells.Item(202, 27).Value = -0.5
$ws.Cells.Item(202, 28).Value = -1
$ws.Cells.Item(202, 29).Value = 0.8999999999999999

# --- Insert new rows 203-207 ---
# Row 203
$ws.Cells.Item(203, 1).Value = 201
$ws.Cells.Item(203, 2).Value = 6769300
$ws.Cells.Item(203, 3).Value = "Croatia HNL"
$ws.Cells.Item(203, 4).Value = "Croatia HNL"
$ws.Cells.Item(203, 5).Value = 45331.54166666666
$ws.Cells.Item(203, 6).Value = "NK Varazdin"
$ws.Cells.Item(203, 7).Value = "NK Rudes"
$ws.Cells.Item(203, 11).Value = 1.6
$ws.Cells.Item(203, 12).Value = 3.8
$ws.Cells.Item(203, 13).Value = 5.5
$ws.Cells.Item(203, 14).Value = 1.6
$ws.Cells.Item(203, 15).Value = 3.8
$ws.Cells.Item(203, 16).Value = 5.75
$ws.Cells.Item(203, 17).Value = -1
$ws.Cells.Item(203, 18).Value = 2.05
$ws.Cells.Item(203, 19).Value = 1.8
$ws.Cells.Item(203, 20).Value = 2.5
$ws.Cells.Item(203, 21).Value = 1.875
$ws.Cells.Item(203, 22).Value = 1.975
$ws.Cells.Item(203, 23).Value = 0
$ws.Cells.Item(203, 24).Value = 0
$ws.Cells.Item(203, 25).Value = 0
$ws.Cells.Item(203, 26).Value = 0
$ws.Cells.Item(203, 27).Value = 0
$idStyleSrc.Copy()
$ws.Cells.Item(203, 1).PasteSpecial(-4122)
$dateStyleSrc.Copy()
$ws.Cells.Item(203, 5).PasteSpecial(-4122)

# Row 204
$ws.Cells.Item(204, 1).Value = 202
$ws.Cells.Item(204, 2).Value = 6788918
$ws.Cells.Item(204, 3).Value = "Croatia HNL"
$ws.Cells.Item(204, 4).Value = "Croatia HNL"
$ws.Cells.Item(204, 5).Value = 45332.45833333334
$ws.Cells.Item(204, 6).Value = "Istra 1961"
$ws.Cells.Item(204, 7).Value = "NK Osijek"
$ws.Cells.Item(204, 11).Value = 3.25
$ws.Cells.Item(204, 12).Value = 3.4
$ws.Cells.Item(204, 13).Value = 2
$ws.Cells.Item(204, 14).Value = 3.1
$ws.Cells.Item(204, 15).Value = 3.1
$ws.Cells.Item(204, 16).Value = 2.15
$ws.Cells.Item(204, 17).Value = 0.25
$ws.Cells.Item(204, 18).Value = 1.975
$ws.Cells.Item(204, 19).Value = 1.875
$ws.Cells.Item(204, 20).Value = 2.25
$ws.Cells.Item(204, 21).Value = 1.825
$ws.Cells.Item(204, 22).Value = 2.025
$ws.Cells.Item(204, 23).Value = 0
$ws.Cells.Item(204, 24).Value = 0
$ws.Cells.Item(204, 25).Value = 0
$ws.Cells.Item(204, 26).Value = 0
$ws.Cells.Item(204, 27).Value = 0
$idStyleSrc.Copy()
$ws.Cells.Item(204, 1).PasteSpecial(-4122)
$dateStyleSrc.Copy()
$ws.Cells.Item(204, 5).PasteSpecial(-4122)

# Row 205
$ws.Cells.Item(205, 1).Value = 203
$ws.Cells.Item(205, 2).Value = 6787891
$ws.Cells.Item(205, 3).Value = "Croatia HNL"
$ws.Cells.Item(205, 4).Value = "Croatia HNL"
$ws.Cells.Item(205, 5).Value = 45332.54861111111
$ws.Cells.Item(205, 6).Value = "NK Lokomotiva Zagreb"
$ws.Cells.Item(205, 7).Value = "Dinamo Zagreb"
$ws.Cells.Item(205, 11).Value = 5.5
$ws.Cells.Item(205, 12).Value = 3.75
$ws.Cells.Item(205, 13).Value = 1.615
$ws.Cells.Item(205, 14).Value = 5.5
$ws.Cells.Item(205, 15).Value = 3.6
$ws.Cells.Item(205, 16).Value = 1.65
$ws.Cells.Item(205, 17).Value = 0.75
$ws.Cells.Item(205, 18).Value = 2.025
$ws.Cells.Item(205, 19).Value = 1.825
$ws.Cells.Item(205, 20).Value = 2.5
$ws.Cells.Item(205, 21).Value = 2.05
$ws.Cells.Item(205, 22).Value = 1.8
$ws.Cells.Item(205, 23).Value = 0
$ws.Cells.Item(205, 24).Value = 0
$ws.Cells.Item(205, 25).Value = 0
$ws.Cells.Item(205, 26).Value = 0
$ws.Cells.Item(205, 27).Value = 0
$idStyleSrc.Copy()
$ws.Cells.Item(205, 1).PasteSpecial(-4122)
$dateStyleSrc.Copy()
$ws.Cells.Item(205, 5).PasteSpecial(-4122)

# Row 206
$ws.Cells.Item(206, 1).Value = 204
$ws.Cells.Item(206, 2).Value = 6788917
$ws.Cells.Item(206, 3).Value = "Croatia HNL"
$ws.Cells.Item(206, 4).Value = "Croatia HNL"
$ws.Cells.Item(206, 5).Value = 45333.45833333334
$ws.Cells.Item(206, 6).Value = "Hajduk Split"
$ws.Cells.Item(206, 7).Value = "Slaven Belupo"
$ws.Cells.Item(206, 11).Value = 1.25
$ws.Cells.Item(206, 12).Value = 5.75
$ws.Cells.Item(206, 13).Value = 9
$ws.Cells.Item(206, 14).Value = 1.25
$ws.Cells.Item(206, 15).Value = 5.75
$ws.Cells.Item(206, 16).Value = 9
$ws.Cells.Item(206, 17).Value = -1.5
$ws.Cells.Item(206, 18).Value = 1.85
$ws.Cells.Item(206, 19).Value = 2
$ws.Cells.Item(206, 20).Value = 2.75
$ws.Cells.Item(206, 21).Value = 1.975
$ws.Cells.Item(206, 22).Value = 1.875
$ws.Cells.Item(206, 23).Value = 0
$ws.Cells.Item(206, 24).Value = 0
$ws.Cells.Item(206, 25).Value = 0
$ws.Cells.Item(206, 26).Value = 0
$ws.Cells.Item(206, 27).Value = 0
$idStyleSrc.Copy()
$ws.Cells.Item(206, 1).PasteSpecial(-4122)
$dateStyleSrc.Copy()
$ws.Cells.Item(206, 5).PasteSpecial(-4122)

# Row 207
$ws.Cells.Item(207, 1).Value = 205
$ws.Cells.Item(207, 2).Value = 6788919
$ws.Cells.Item(207, 3).Value = "Croatia HNL"
$ws.Cells.Item(207, 4).Value = "Croatia HNL"
$ws.Cells.Item(207, 5).Value = 45333.54861111111
$ws.Cells.Item(207, 6).Value = "HNK Gorica"
$ws.Cells.Item(207, 7).Value = "HNK Rijeka"
$ws.Cells.Item(207, 11).Value = 4.75
$ws.Cells.Item(207, 12).Value = 3.5
$ws.Cells.Item(207, 13).Value = 1.727
$ws.Cells.Item(207, 14).Value = 5
$ws.Cells.Item(207, 15).Value = 3.6
$ws.Cells.Item(207, 16).Value = 1.666
$ws.Cells.Item(207, 17).Value = 0.75
$ws.Cells.Item(207, 18).Value = 1.95
$ws.Cells.Item(207, 19).Value = 1.9
$ws.Cells.Item(207, 20).Value = 2.5
$ws.Cells.Item(207, 21).Value = 1.975
$ws.Cells.Item(207, 22).Value = 1.875
$ws.Cells.Item(207, 23).Value = 0
$ws.Cells.Item(207, 24).Value = 0
$ws.Cells.Item(207, 25).Value = 0
$ws.Cells.Item(207, 26).Value = 0
$ws.Cells.Item(207, 27).Value = 0
$idStyleSrc.Copy()
$ws.Cells.Item(207, 1).PasteSpecial(-4122)
$dateStyleSrc.Copy()
$ws.Cells.Item(207, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
